$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 & 4: reword the "main menu option" test case descriptions and narrow
# their test values down to single options instead of option ranges.
$ws.Range("C3").Value = "To test if the program accepts a valid input for main menu selection."
$ws.Range("C4").Value = "To test if the program triggers validation for an invalid input for main menu selection."
$ws.Range("D3").Value = "Input: `nOption - 0"
$ws.Range("B3").Value = "Test valid main menu option is accepted"
$ws.Range("D4").Value = "Input: `nOption - 5"

# Update the active selection to D5 (no explicit frozen/scrolled top-left cell anymore).
$ws.Range("D5").Select()
